# ExcelDemo / JsonDemo.xlsx update
# - Drops the old "ProjectID" column, shifting ProjectName into column A
# - Replaces the generic Key1/Value1/Key2/Value2 placeholder columns with
#   real task-tracking columns (Environment/Task/Manager/Author/Target Team)
# - Replaces the sample rows with the new project-task data
# - Updates column widths, active-cell selection and page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "ProjectName"
$ws.Range("B1").Value = "PhaseID"
$ws.Range("C1").Value = "VerticalID"
$ws.Range("D1").Value = "Environment "
$ws.Range("E1").Value = "Task"
$ws.Range("F1").Value = "Manager "
$ws.Range("G1").Value = "Author"
$ws.Range("H1").Value = "Target Team"

# ---- Row 2 ----
$ws.Range("A2").Value = "Excel Security"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = "Laura's Test Environment"
$ws.Range("E2").Value = "Add detail page"
$ws.Range("F2").Value = "Jimmy"
$ws.Range("G2").Value = "Jimmy"
$ws.Range("H2").Value = "UX"

# ---- Row 3 ----
$ws.Range("A3").Value = "Excel Security"
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = "Laura's Test Environment"
$ws.Range("E3").Value = "Convert to UpdatePackage"
$ws.Range("F3").Value = "Martha"
$ws.Range("G3").Value = "Same"
$ws.Range("H3").Value = "DB"

# ---- Row 4 ----
$ws.Range("A4").Value = "Excel Security"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = "Laura's Build Environment"
$ws.Range("E4").Value = "Update Contact page"
$ws.Range("F4").Value = "Bocephus"
$ws.Range("G4").Value = "Tooty"
$ws.Range("H4").Value = "Back End"

# ---- Row 5 ----
$ws.Range("A5").Value = "Excel Security"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = "BR549"
$ws.Range("E5").Value = "Test Register function"
$ws.Range("F5").Value = "Marty"
$ws.Range("G5").Value = "Sarah"
$ws.Range("H5").Value = "Architecture"

# ---- Column widths (re-fit to the new, narrower column A and the new
#      data so every column matches its longest entry again) ----
$ws.Columns("A").ColumnWidth = 16.8333
$ws.Columns("B").ColumnWidth = 7.3333
$ws.Columns("C").ColumnWidth = 8.8333
$ws.Columns("D").ColumnWidth = 18.8333
$ws.Columns("E").ColumnWidth = 20.5
$ws.Columns("F").ColumnWidth = 14.8333
$ws.Columns("G").ColumnWidth = 13.8333
$ws.Columns("H").ColumnWidth = 15.6667

# ---- Page setup: explicit portrait orientation ----
$ws.PageSetup.Orientation = 1

# ---- Selection moved one row down from where it used to be ----
$ws.Range("H6").Select()
